$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blue")
$ws.Activate()

# Project start date shifts back by one day (literal value, not a formula)
$ws.Range("C5").Value = $ws.Range("C5").Value - 1

# Scroll/selection state as left by the author
$ws.Range("M11").Select()
$excel.ActiveWindow.ScrollRow = 4
